$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header cells: "JAVA" -> "Relational Calculus", "Operating Systems" -> "Algebra"
$ws.Range("C1").Value = "Relational Calculus"
$ws.Range("D1").Value = "Algebra"

# Widen column C to fit the new longer header text
$ws.Columns("C").ColumnWidth = 19.65

# Move the active selection to D1
$ws.Range("D1").Select()
